# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 2068
    4  = 129
    5  = 43
    7  = 1707
    9  = 704
    11 = 105
    13 = 103
    14 = 230
    16 = 141
    19 = 3950
    22 = 446
    23 = 378
    24 = 839
    25 = 615
    27 = 9
    29 = 1760
    30 = 28
    31 = 31
    32 = 59
    33 = 176
    34 = 16
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
